$d = $word.ActiveDocument

# --- Step 1: locate the "Indirect Recursion...- Cuando..." paragraph (last paragraph) ---
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range

# Sanity check we have the right paragraph.
if ($r.Text -notlike "Indirect Recursion*") {
    throw "Unexpected last paragraph: $($r.Text)"
}

# --- Step 2: remove the hidden _GoBack bookmark sitting at the end of this paragraph ---
$bm = $r.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Step 3: fix the run split around ".- " -> the hyphen moves into the previous run ---
# Locate "Recursion.- Cuando" inside the paragraph to compute exact character offsets.
$search = $d.Range($r.Start, $r.End)
$found = $search.Find.Execute("Recursion.- Cuando")
if (-not $found) {
    throw "Could not find 'Recursion.- Cuando' text"
}

$dotStart = $search.Start + 9
$dotEnd = $dotStart + 1
$hyStart = $dotEnd
$hyEnd = $hyStart + 2

# Edit right-to-left so earlier offsets stay valid.
$hyRange = $d.Range($hyStart, $hyEnd)
if ($hyRange.Text -ne "- ") {
    throw "Unexpected hyphen run text: [$($hyRange.Text)]"
}
$hyRange.Text = " "

$dotRange = $d.Range($dotStart, $dotEnd)
if ($dotRange.Text -ne ".") {
    throw "Unexpected dot run text: [$($dotRange.Text)]"
}
$dotRange.Text = ".-"

# --- Step 4: insert the new paragraphs right after this paragraph ---
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$insertionPoint = $d.Range($r.End, $r.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = @"
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Suma de números naturales con </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>recursión.-</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Formula incluye recursión:</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>/</w:t></w:r><w:r><w:tab/><w:t>0</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>n = 0</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>O(n)</w:t></w:r><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>calls</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: n + 1</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Sum(n)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>\</w:t></w:r><w:r><w:tab/><w:t>Sum(n-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1)+</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>n</w:t></w:r><w:r><w:tab/><w:t>n &gt; 0</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Las definiciones matemáticas con recursión son fácilmente convertibles a funciones recursivas.</w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Se puede lograr una formula para evitar la recursión. En este </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>caso :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">  n (n+1) / 2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
<w:p $ns><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Esto hará que la complejidad de tiempo sea </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>O(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>1)</w:t></w:r></w:p>
"@

$insertionPoint.InsertXML($newXml)

Write-Host "Edit applied. Paragraph count: $($d.Paragraphs.Count)"
